$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.363.27'
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').Value = '2.587.08'
$ws.Range('E3').Value = '  -2.46%  '
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').Value = '560.35'
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('D6').Value = '143.34'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('D9').Value = '2.593.98'
$ws.Range('E9').Value = '  -3.19%  '
$ws.Range('D10').Value = '6.65'
$ws.Range('E10').Value = '  -3.09%  '
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('E12').Value = '  +10.78%  '
$ws.Range('D13').Value = '0.358'
$ws.Range('E13').Value = '  +4.54%  '
$ws.Range('D14').Value = '3.043.60'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '23.36'
$ws.Range('E15').Value = '  +6.70%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '59.303.97'
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '2.573.91'
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('D19').Value = '4.60'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = '337.82'
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').Value = '10.38'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '64.05'
$ws.Range('E24').Value = '  -4.30%  '
$ws.Range('E25').Value = '  +5.48%  '
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('D28').Value = '7.40'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').Value = '0.0₃0778'
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  -3.17%  '
$ws.Range('D33').Value = '159.06'
$ws.Range('E33').Value = '  +2.89%  '
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('D37').Value = '0.884'
$ws.Range('E37').Value = '  -3.33%  '
$ws.Range('D38').Value = '0.875'
$ws.Range('E38').Value = '  -3.75%  '
$ws.Range('D39').Value = '37.46'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').Value = '3.69'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('D42').Value = '293.33'
$ws.Range('E42').Value = '  -4.14%  '
$ws.Range('D43').Value = '133.41'
$ws.Range('E43').Value = '  +5.04%  '
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').Value = '0.598'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('D47').Value = '0.0537'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').Value = '10.64'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('D49').Value = '0.0235'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('D50').Value = '18.75'
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('D51').Value = '1.954.78'
$ws.Range('E51').Value = '  -0.77%  '
